# Update Olympic's fixtures list
# Mark "Data availability" (column G) as "Y" for the rounds that were
# previously marked "N" now that the underlying match data has become
# available, and move the selection to reflect the last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixtures")

# Rows 12,13 and 16-20 flip from "N" to "Y" in the Data availability column (G)
$ws.Range("G12").Value = "Y"
$ws.Range("G13").Value = "Y"
$ws.Range("G16").Value = "Y"
$ws.Range("G17").Value = "Y"
$ws.Range("G18").Value = "Y"
$ws.Range("G19").Value = "Y"
$ws.Range("G20").Value = "Y"

# Leave the active selection on the last-touched cell
$ws.Range("G20").Select()
